$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Variables")

# Remove the "row_id" row (row 2) - rows below shift up.
$ws.Rows.Item(2).Delete()

# Leave the selection where Excel would naturally land after the delete.
$ws.Activate()
$ws.Range("H12").Select()
